# AI算法与应用-2024工作进展.xlsx
# 毛彬 added his weekly progress entry (new row 16) on sheet "001毛彬".
# Row 15 picks up "plain" default formatting in the process (matches
# the other data rows A2:D14 - default style for col A, wrap-text style
# for cols B:D) instead of its previous, row-specific duplicate styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("001毛彬")
$ws.Select()

# --- normalize row 15's formatting to the same styles used by every
#     other row in the table (default style for A, wrap-text style for
#     B:D) by copying the formats down from existing "normal" cells ---
$ws.Range("A1").Copy()
$ws.Range("A15").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B2").Copy()
$ws.Range("B15:D15").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- append the new weekly entry as row 16 ---
$ws.Range("A16").Value = "20240430-20240506"
$ws.Range("B16").Value = "1.数据通过CWT转换为图片，输入到swimtransform、mobilenetv2、resnet50模型中，模型出现严重过拟合现象`n2.补第5、6次两周六道题"
$ws.Range("C16").Value = "调了好几次超参数，依旧欠拟合严重，没有头绪。估计是数据转换方法需要换"
$ws.Range("D16").Value = "1.查找文献，换一种数据转换的方法，进行模型的训练`n2.补第7次两周六道题"

# new row inherits the same wrap-text formatting as the rest of the table
$ws.Range("B2").Copy()
$ws.Range("B16:D16").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- view state: user scrolled down and zoomed in, with G15 selected ---
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 1
$win.Zoom = 115

$ws.Range("G15").Select()
